# Update BOC USD rates (auto)
# Adds the newly published 2026-01-04 10:30:00 rate row to "All Published
# Values", grows the sheet's used range / AutoFilter accordingly, and bumps
# the "Daily Summary" publish count for 2026-01-04 from 3 to 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "All Published Values": append row 27 with the new published value.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("All Published Values")

$newRow = $ws1.Range("A27:J27")
# Force text entry so date-/number-looking strings aren't auto-converted
# (the sheet stores every data value as text, matching the rest of the
# column).
$newRow.NumberFormat = "@"

$ws1.Range("A27").Value = "2026-01-04"
$ws1.Range("B27").Value = "2026-01-04 10:30:00"
$ws1.Range("C27").Value = "697.35"
$ws1.Range("D27").Value = "697.35"
$ws1.Range("E27").Value = "700.29"
$ws1.Range("F27").Value = "700.29"
$ws1.Range("G27").Value = "702.88"
$ws1.Range("H27").Value = "2026/01/04 10:30:00"
$ws1.Range("I27").Value = "2026-01-04 03:49:03"
$ws1.Range("J27").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Drop back to the default (unstyled) look used by the rest of the table.
$newRow.Style = "Normal"

# Re-apply the AutoFilter so its range grows from A1:J26 to A1:J27 (toggle
# off, then back on, since the range already has a filter applied).
$ws1.Range("A1:J27").AutoFilter()
$ws1.Range("A1:J27").AutoFilter()

# The hidden _xlnm._FilterDatabase defined name also needs to track the
# new range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -eq "All Published Values!_FilterDatabase") {
        $nm.RefersTo = "='All Published Values'!`$A`$1:`$J`$27"
    }
}

# ---------------------------------------------------------------------
# 2. "Daily Summary": the 2026-01-04 row now has 4 publishes instead of 3.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("B6").Value = 4
